# Adds the "Start Date", "Acceptance Date" and "Salary (Annual)" columns to
# the candidate details sheet, renames the "S.No" header to "S No", renames
# the worksheet, and tweaks a couple of view/page-setup details.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---------------------------------------------
$ws.Name = "Candidates Data"

# --- New data (per-row: Start Date, Acceptance Date, Salary (Annual)) ----
# Dates are stored as Excel serial day numbers (1899-12-30 epoch).
$data = @(
    @(45992, 45981, 650000),
    @(45996, 45983, 580000),
    @(45994, 45982, 520000),
    @(46001, 45986, 600000),
    @(45998, 45984, 570000),
    @(46003, 45987, 620000),
    @(45999, 45985, 500000),
    @(46006, 45989, 640000),
    @(46009, 45991, 720000),
    @(46011, 45992, 800000)
)

# --- Headers for the three new columns (written before the "S No" header
#     so the shared-string table is rebuilt in this same order) -----------
$ws.Range("G1").Value = "Start Date"
$ws.Range("H1").Value = "Acceptance Date"
$ws.Range("I1").Value = "Salary (Annual)"

$ws.Range("G1:I1").Font.Bold = $true
$ws.Range("G1:I1").HorizontalAlignment = -4131

# --- Row data for the new columns -----------------------------------------
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i

    $ws.Cells.Item($r, 7).Value = $data[$i][0]
    $ws.Cells.Item($r, 8).Value = $data[$i][1]
    $ws.Cells.Item($r, 9).Value = $data[$i][2]
}

$ws.Range("G2:H11").NumberFormat = "d-mmm-yy"
$ws.Range("G2:I11").HorizontalAlignment = -4131

# --- Column widths for the new columns ------------------------------------
$ws.Columns.Item(7).ColumnWidth = 16.27
$ws.Columns.Item(8).ColumnWidth = 21.27
$ws.Columns.Item(9).ColumnWidth = 14

# --- Rename "S.No" -> "S No" (written last so it lands at the end of the
#     rebuilt shared-string table) -----------------------------------------
$ws.Range("A1").Value = "S No"

# --- View tweaks ------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("C15").Select()

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
